# Update the TAS2 "survey" sheet constraint text so that the allowed
# EU code range changes from 9-1000 to 9-100, and the allowed respondent
# age range changes from 5-7 to 5-12 years (English + French messages).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# p_recorderID constraint messages (row 2): "...between 9 and 1000" -> "...between 9 and 100"
$ws.Range("I2").Value = "The code must be a two-digit number between 9 and 100"
$ws.Range("J2").Value = "Le code doit être un nombre à deux chiffres entre 9 et 100"

# p_age_yrs constraint + messages (row 7): 5-7 years -> 5-12 years
$ws.Range("H7").Value = ". >= 5 and . <= 12"
$ws.Range("I7").Value = "The age must be between 5 and 12 years"
$ws.Range("J7").Value = "L'age doit être compris entre 5 et 12 ans"
